$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title textbox ("TextBox 6"): "6.042 Microquiz April 10, 2013" (32pt)
#     -> "6.042 Microquiz9w" (40pt), moved/resized to its new autofit box. ---
$title = $s.Shapes.Item("TextBox 6")

$tr = $title.TextFrame.TextRange
$tr.Text = "6.042 Microquiz9w"
$tr.Font.Size = 40

# Target EMU values (from the authoritative OOXML): off (2228349,318039),
# ext (4711646,707886). PowerPoint's COM surface works in points (1 pt =
# 12700 EMU); the values below are chosen so the point -> EMU conversion
# lands exactly on those targets.
$title.Left = 175.4605905511811
$title.Top = 25.042480314960628
$title.Width = 370.99578740157483
$title.Height = 55.739094488188975
